$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.973.08'
$ws.Range("E2").Value = '  +3.40%  '

$ws.Range("D3").Value = '3.407.71'
$ws.Range("E3").Value = '  +1.84%  '

$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = '577.73'
$ws.Range("E5").Value = '  +2.36%  '

$ws.Range("D6").Value = '137.67'
$ws.Range("E6").Value = '  +5.47%  '

$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("E8").Value = '  +1.01%  '

$ws.Range("E9").Value = '  +0.91%  '

$ws.Range("E10").Value = '  +7.34%  '

$ws.Range("E11").Value = '  +4.67%  '

$ws.Range("D12").Value = '3.992.08'
$ws.Range("E12").Value = '  +2.02%  '

$ws.Range("E13").Value = '  +2.34%  '

$ws.Range("E14").Value = '  +5.89%  '

$ws.Range("D15").Value = '3.405.94'
$ws.Range("E15").Value = '  +1.84%  '

$ws.Range("D16").Value = '25.47'
$ws.Range("E16").Value = '  +3.51%  '

$ws.Range("D17").Value = '61.971.93'
$ws.Range("E17").Value = '  +3.02%  '

$ws.Range("E18").Value = '  +6.11%  '

$ws.Range("D19").Value = '5.90'
$ws.Range("E19").Value = '  +3.96%  '

$ws.Range("D20").Value = '9.50'
$ws.Range("E20").Value = '  +5.59%  '

$ws.Range("D21").Value = '389.43'
$ws.Range("E21").Value = '  +9.93%  '

$ws.Range("E22").Value = '  +2.57%  '

$ws.Range("D23").Value = '3.544.21'
$ws.Range("E23").Value = '  +1.94%  '

$ws.Range("E24").Value = '  +14.73%  '

$ws.Range("E25").Value = '  -0.05%  '

$ws.Range("D26").Value = '71.55'
$ws.Range("E26").Value = '  +3.29%  '

$ws.Range("D27").Value = '7.70'
$ws.Range("E27").Value = '  +3.14%  '

$ws.Range("E28").Value = '  -0.06%  '

$ws.Range("E29").Value = '  +0.23%  '

$ws.Range("D30").Value = '8.31'
$ws.Range("E30").Value = '  +5.17%  '

$ws.Range("E31").Value = '  +4.38%  '

$ws.Range("E32").Value = '  +2.59%  '

$ws.Range("B33").Value = 'RenzoRestakedETH'
$ws.Range("C33").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D33").Value = '3.438.81'
$ws.Range("E33").Value = '  +1.84%  '

$ws.Range("B34").Value = 'USDe'
$ws.Range("C34").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  +0.01%  '

$ws.Range("E35").Value = '  +3.00%  '

$ws.Range("E36").Value = '  +1.79%  '

$ws.Range("D37").Value = '7.00'
$ws.Range("E37").Value = '  +2.22%  '

$ws.Range("E38").Value = '  +4.13%  '

$ws.Range("D39").Value = '164.43'
$ws.Range("E39").Value = '  +3.87%  '

$ws.Range("E40").Value = '  +3.56%  '

$ws.Range("D41").Value = '1.79'
$ws.Range("E41").Value = '  +13.54%  '

$ws.Range("E42").Value = '  +5.29%  '

$ws.Range("B43").Value = 'ONDO'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D43").Value = '1.23'
$ws.Range("E43").Value = '  +2.88%  '

$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").Value = '0.999'
$ws.Range("E44").Value = '  -0.01%  '

$ws.Range("E45").Value = '  +2.23%  '

$ws.Range("D46").Value = '25.04'
$ws.Range("E46").Value = '  +6.38%  '

$ws.Range("D47").Value = '41.67'
$ws.Range("E47").Value = '  +1.97%  '

$ws.Range("E48").Value = '  +2.19%  '

$ws.Range("E49").Value = '  +3.82%  '

$ws.Range("D50").Value = '2.376.93'
$ws.Range("E50").Value = '  +9.44%  '

$ws.Range("E51").Value = '  +6.66%  '
